$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.764753333333333
$ws.Range("H2").Value = 14.29426
$ws.Range("I2").Value = 0.2966169987831952
$ws.Range("J2").Value = 0.2966169987831952
$ws.Range("M2").Value = 7.236245333333333
$ws.Range("N2").Value = 21.708736
$ws.Range("O2").Value = 0.6630432242917509
$ws.Range("P2").Value = 0.6630432242917509
$ws.Range("Q2").Value = 34.47892407281778
$ws.Range("R2").Value = 310.3103166553599
$ws.Range("S2").Value = 0.1966698912529521
$ws.Range("T2").Value = 0.1966698912529521

# Row 3
$ws.Range("G3").Value = 4.764753333333333
$ws.Range("H3").Value = 14.29426
$ws.Range("I3").Value = 0.2966169987831952
$ws.Range("J3").Value = 0.2966169987831952
$ws.Range("O3").Value = 0.003787629702975075
$ws.Range("P3").Value = 0.003787629702975075
$ws.Range("Q3").Value = 0.19696060854
$ws.Range("R3").Value = 1.77264547686
$ws.Range("S3").Value = 0.001123475354998552
$ws.Range("T3").Value = 0.001123475354998552

# Row 4
$ws.Range("G4").Value = 4.764753333333333
$ws.Range("H4").Value = 14.29426
$ws.Range("I4").Value = 0.2966169987831952
$ws.Range("J4").Value = 0.2966169987831952
$ws.Range("M4").Value = 3.636103333333333
$ws.Range("N4").Value = 10.90831
$ws.Range("O4").Value = 0.3331691460052741
$ws.Range("P4").Value = 0.3331691460052741
$ws.Range("Q4").Value = 17.32513547784444
$ws.Range("R4").Value = 155.9262193006
$ws.Range("S4").Value = 0.09882363217524456
$ws.Range("T4").Value = 0.09882363217524456

# Row 5
$ws.Range("I5").Value = 0.5337607564504776
$ws.Range("J5").Value = 0.5337607564504775
$ws.Range("M5").Value = 7.236245333333333
$ws.Range("N5").Value = 21.708736
$ws.Range("O5").Value = 0.6630432242917509
$ws.Range("P5").Value = 0.6630432242917509
$ws.Range("Q5").Value = 62.04464568855467
$ws.Range("R5").Value = 558.401811196992
$ws.Range("S5").Value = 0.3539064529573286
$ws.Range("T5").Value = 0.3539064529573286

# Row 6
$ws.Range("I6").Value = 0.5337607564504776
$ws.Range("J6").Value = 0.5337607564504775
$ws.Range("O6").Value = 0.003787629702975075
$ws.Range("P6").Value = 0.003787629702975075
$ws.Range("S6").Value = 0.002021688095414274
$ws.Range("T6").Value = 0.002021688095414273

# Row 7
$ws.Range("I7").Value = 0.5337607564504776
$ws.Range("J7").Value = 0.5337607564504775
$ws.Range("M7").Value = 3.636103333333333
$ws.Range("N7").Value = 10.90831
$ws.Range("O7").Value = 0.3331691460052741
$ws.Range("P7").Value = 0.3331691460052741
$ws.Range("Q7").Value = 31.17649175939667
$ws.Range("R7").Value = 280.58842583457
$ws.Range("S7").Value = 0.1778326153977347
$ws.Range("T7").Value = 0.1778326153977347

# Row 8
$ws.Range("G8").Value = 2.724753333333334
$ws.Range("H8").Value = 8.17426
$ws.Range("I8").Value = 0.1696222447663273
$ws.Range("J8").Value = 0.1696222447663273
$ws.Range("M8").Value = 7.236245333333333
$ws.Range("N8").Value = 21.708736
$ws.Range("O8").Value = 0.6630432242917509
$ws.Range("P8").Value = 0.6630432242917509
$ws.Range("Q8").Value = 19.71698359281778
$ws.Range("R8").Value = 177.45285233536
$ws.Range("S8").Value = 0.1124668800814702
$ws.Range("T8").Value = 0.1124668800814702

# Row 9
$ws.Range("G9").Value = 2.724753333333334
$ws.Range("H9").Value = 8.17426
$ws.Range("I9").Value = 0.1696222447663273
$ws.Range("J9").Value = 0.1696222447663273
$ws.Range("O9").Value = 0.003787629702975075
$ws.Range("P9").Value = 0.003787629702975075
$ws.Range("Q9").Value = 0.11263312854
$ws.Range("R9").Value = 1.01369815686
$ws.Range("S9").Value = 0.0006424662525622497
$ws.Range("T9").Value = 0.0006424662525622496

# Row 10
$ws.Range("G10").Value = 2.724753333333334
$ws.Range("H10").Value = 8.17426
$ws.Range("I10").Value = 0.1696222447663273
$ws.Range("J10").Value = 0.1696222447663273
$ws.Range("M10").Value = 3.636103333333333
$ws.Range("N10").Value = 10.90831
$ws.Range("O10").Value = 0.3331691460052741
$ws.Range("P10").Value = 0.3331691460052741
$ws.Range("Q10").Value = 9.907484677844446
$ws.Range("R10").Value = 89.1673621006
$ws.Range("S10").Value = 0.05651289843229483
$ws.Range("T10").Value = 0.05651289843229482
